# Update crypto prices/volume percentages per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '38.315.74'
$cell.Style = "Normal"

$cell = $ws.Range('E2')
$cell.NumberFormat = "@"
$cell.Value = '  +1.46%  '
$cell.Style = "Normal"

$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '2.097.08'
$cell.Style = "Normal"

$cell = $ws.Range('E3')
$cell.NumberFormat = "@"
$cell.Value = '  +3.17%  '
$cell.Style = "Normal"

$cell = $ws.Range('E4')
$cell.NumberFormat = "@"
$cell.Value = '  -0.15%  '
$cell.Style = "Normal"

$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '228.78'
$cell.Style = "Normal"

$cell = $ws.Range('E5')
$cell.NumberFormat = "@"
$cell.Value = '  +0.55%  '
$cell.Style = "Normal"

$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '0.614'
$cell.Style = "Normal"

$cell = $ws.Range('E6')
$cell.NumberFormat = "@"
$cell.Value = '  +0.97%  '
$cell.Style = "Normal"

$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '61.12'
$cell.Style = "Normal"

$cell = $ws.Range('E7')
$cell.NumberFormat = "@"
$cell.Value = '  +1.71%  '
$cell.Style = "Normal"

$cell = $ws.Range('E8')
$cell.NumberFormat = "@"
$cell.Value = '  -0.03%  '
$cell.Style = "Normal"

$cell = $ws.Range('E9')
$cell.NumberFormat = "@"
$cell.Value = '  +0.83%  '
$cell.Style = "Normal"

$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '0.0848'
$cell.Style = "Normal"

$cell = $ws.Range('E10')
$cell.NumberFormat = "@"
$cell.Value = '  +3.05%  '
$cell.Style = "Normal"

$cell = $ws.Range('E11')
$cell.NumberFormat = "@"
$cell.Value = '  +0.41%  '
$cell.Style = "Normal"

$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '2.407.01'
$cell.Style = "Normal"

$cell = $ws.Range('E12')
$cell.NumberFormat = "@"
$cell.Value = '  +3.03%  '
$cell.Style = "Normal"

$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '14.79'
$cell.Style = "Normal"

$cell = $ws.Range('E13')
$cell.NumberFormat = "@"
$cell.Value = '  +2.45%  '
$cell.Style = "Normal"

$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '22.36'
$cell.Style = "Normal"

$cell = $ws.Range('E14')
$cell.NumberFormat = "@"
$cell.Value = '  +6.51%  '
$cell.Style = "Normal"

$cell = $ws.Range('B15')
$cell.NumberFormat = "@"
$cell.Value = 'Polygon'
$cell.Style = "Normal"

$cell = $ws.Range('C15')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$cell.Style = "Normal"

$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '0.778'
$cell.Style = "Normal"

$cell = $ws.Range('E15')
$cell.NumberFormat = "@"
$cell.Value = '  +2.69%  '
$cell.Style = "Normal"

$cell = $ws.Range('B16')
$cell.NumberFormat = "@"
$cell.Value = 'Polkadot'
$cell.Style = "Normal"

$cell = $ws.Range('C16')
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell.Style = "Normal"

$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '5.46'
$cell.Style = "Normal"

$cell = $ws.Range('E16')
$cell.NumberFormat = "@"
$cell.Value = '  +5.58%  '
$cell.Style = "Normal"

$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '2.098.16'
$cell.Style = "Normal"

$cell = $ws.Range('E17')
$cell.NumberFormat = "@"
$cell.Value = '  +3.28%  '
$cell.Style = "Normal"

$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '38.268.71'
$cell.Style = "Normal"

$cell = $ws.Range('E18')
$cell.NumberFormat = "@"
$cell.Value = '  +1.46%  '
$cell.Style = "Normal"

$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '70.45'
$cell.Style = "Normal"

$cell = $ws.Range('E19')
$cell.NumberFormat = "@"
$cell.Value = '  +1.23%  '
$cell.Style = "Normal"

$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '6.02'
$cell.Style = "Normal"

$cell = $ws.Range('E20')
$cell.NumberFormat = "@"
$cell.Value = '  +1.51%  '
$cell.Style = "Normal"

$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0834'
$cell.Style = "Normal"

$cell = $ws.Range('E21')
$cell.NumberFormat = "@"
$cell.Value = '  +1.47%  '
$cell.Style = "Normal"

$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '225.00'
$cell.Style = "Normal"

$cell = $ws.Range('E22')
$cell.NumberFormat = "@"
$cell.Value = '  +0.66%  '
$cell.Style = "Normal"

$cell = $ws.Range('E23')
$cell.NumberFormat = "@"
$cell.Value = '  +0.03%  '
$cell.Style = "Normal"

$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '2.40'
$cell.Style = "Normal"

$cell = $ws.Range('E24')
$cell.NumberFormat = "@"
$cell.Value = '  -0.79%  '
$cell.Style = "Normal"

$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '2.32'
$cell.Style = "Normal"

$cell = $ws.Range('E25')
$cell.NumberFormat = "@"
$cell.Value = '  +3.49%  '
$cell.Style = "Normal"

$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '169.84'
$cell.Style = "Normal"

$cell = $ws.Range('E26')
$cell.NumberFormat = "@"
$cell.Value = '  +1.57%  '
$cell.Style = "Normal"

$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '9.39'
$cell.Style = "Normal"

$cell = $ws.Range('E27')
$cell.NumberFormat = "@"
$cell.Value = '  +1.50%  '
$cell.Style = "Normal"

$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '0.131'
$cell.Style = "Normal"

$cell = $ws.Range('E28')
$cell.NumberFormat = "@"
$cell.Value = '  +0.86%  '
$cell.Style = "Normal"

$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '19.05'
$cell.Style = "Normal"

$cell = $ws.Range('E29')
$cell.NumberFormat = "@"
$cell.Value = '  +1.41%  '
$cell.Style = "Normal"

$cell = $ws.Range('E30')
$cell.NumberFormat = "@"
$cell.Value = '  +7.57%  '
$cell.Style = "Normal"

$cell = $ws.Range('E31')
$cell.NumberFormat = "@"
$cell.Value = '  -0.27%  '
$cell.Style = "Normal"

$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '2.34'
$cell.Style = "Normal"

$cell = $ws.Range('E32')
$cell.NumberFormat = "@"
$cell.Value = '  +7.31%  '
$cell.Style = "Normal"

$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '4.75'
$cell.Style = "Normal"

$cell = $ws.Range('E33')
$cell.NumberFormat = "@"
$cell.Value = '  +6.21%  '
$cell.Style = "Normal"

$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '4.44'
$cell.Style = "Normal"

$cell = $ws.Range('E34')
$cell.NumberFormat = "@"
$cell.Value = '  +1.25%  '
$cell.Style = "Normal"

$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '0.0605'
$cell.Style = "Normal"

$cell = $ws.Range('E35')
$cell.NumberFormat = "@"
$cell.Value = '  +0.38%  '
$cell.Style = "Normal"

$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '2.39'
$cell.Style = "Normal"

$cell = $ws.Range('E36')
$cell.NumberFormat = "@"
$cell.Value = '  +5.26%  '
$cell.Style = "Normal"

$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '6.42'
$cell.Style = "Normal"

$cell = $ws.Range('E37')
$cell.NumberFormat = "@"
$cell.Value = '  +1.73%  '
$cell.Style = "Normal"

$cell = $ws.Range('E38')
$cell.NumberFormat = "@"
$cell.Value = '  +6.37%  '
$cell.Style = "Normal"

$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"

$cell = $ws.Range('E39')
$cell.NumberFormat = "@"
$cell.Value = '  -0.10%  '
$cell.Style = "Normal"

$cell = $ws.Range('E40')
$cell.NumberFormat = "@"
$cell.Value = '  +1.89%  '
$cell.Style = "Normal"

$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '1.549.56'
$cell.Style = "Normal"

$cell = $ws.Range('E41')
$cell.NumberFormat = "@"
$cell.Value = '  +0.75%  '
$cell.Style = "Normal"

$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '100.07'
$cell.Style = "Normal"

$cell = $ws.Range('E42')
$cell.NumberFormat = "@"
$cell.Value = '  +4.55%  '
$cell.Style = "Normal"

$cell = $ws.Range('E43')
$cell.NumberFormat = "@"
$cell.Value = '  +1.74%  '
$cell.Style = "Normal"

$cell = $ws.Range('E44')
$cell.NumberFormat = "@"
$cell.Value = '  +1.04%  '
$cell.Style = "Normal"

$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '0.0909'
$cell.Style = "Normal"

$cell = $ws.Range('E45')
$cell.NumberFormat = "@"
$cell.Value = '  -0.24%  '
$cell.Style = "Normal"

$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '4.13'
$cell.Style = "Normal"

$cell = $ws.Range('E46')
$cell.NumberFormat = "@"
$cell.Value = '  -1.41%  '
$cell.Style = "Normal"

$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '1.12'
$cell.Style = "Normal"

$cell = $ws.Range('E47')
$cell.NumberFormat = "@"
$cell.Value = '  +1.08%  '
$cell.Style = "Normal"

$cell = $ws.Range('E48')
$cell.NumberFormat = "@"
$cell.Value = '  +5.35%  '
$cell.Style = "Normal"

$cell = $ws.Range('E49')
$cell.NumberFormat = "@"
$cell.Value = '  +2.55%  '
$cell.Style = "Normal"

$cell = $ws.Range('E50')
$cell.NumberFormat = "@"
$cell.Value = '  +0.90%  '
$cell.Style = "Normal"

$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '2.292.40'
$cell.Style = "Normal"

$cell = $ws.Range('E51')
$cell.NumberFormat = "@"
$cell.Value = '  +2.97%  '
$cell.Style = "Normal"
